$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - force text format first so COM does not
# auto-coerce numeric-looking strings (e.g. "1.00" -> 1) on assignment
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.684.56"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.094.81"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.80"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.30"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0843"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.398.92"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.87"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.37"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.081.36"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.583.64"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.09"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.04"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.38"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.44"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.46"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.14"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.77"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.49"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0607"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.56"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.56"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.55"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.546.38"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.21"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0221"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0917"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.17"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.04"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.289.94"

# Volume 1h % (column E) updates - same text-format guard
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.37%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.61%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.34%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.17%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.25%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.95%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.09%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.49%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.59%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.20%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +9.18%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.84%  "
